# Add a new "2022-Q4" quarter: insert a new worksheet with its fund
# holdings right after "总计" (pushing 2022-Q3 and everything after it
# one slot later), and add the matching summary row to "总计".

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q4 and bump the other
#    rows' running index (column A) down by one.
# ---------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 1.24

for ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------
# 2) New "2022-Q4" worksheet (same layout as the other quarters),
#    inserted right after "总计" / right before "2022-Q3".
# ---------------------------------------------------------------
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$q4.Name = "2022-Q4"

$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q3.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

$q4.Range("B2:G5").NumberFormat = "@"

$rows = @(
    @(0, "100032", "富国中证红利指数增强A",     "59.06", "91.28", "1.85", "1.0926", 6),
    @(1, "008682", "富国中证红利指数增强C",     "5.45",  "91.28", "1.85", "0.1008", 6),
    @(2, "515300", "嘉实沪深300红利低波动ETF", "1.08",  "99.12", "2.96", "0.0320", 9),
    @(3, "510290", "南方上证380ETF",           "1.65",  "99.21", "0.88", "0.0145", 6)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------
# 3) Restore the previously-active "2020-Q4" tab as the selected
#    sheet (Worksheets.Add() above made the new sheet the active one).
# ---------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
